# Apply the F-column ("想去人数" / wanted-to-go count) updates described by the diff.
# The same underlying events are listed on multiple sheets (展览, 演出, and the
# aggregated 全部类型 sheet), so each one needs to be updated in place.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 675
$wsExhibition.Range("F4").Value = 934
$wsExhibition.Range("F5").Value = 697
$wsExhibition.Range("F12").Value = 371
$wsExhibition.Range("F15").Value = 332
$wsExhibition.Range("F16").Value = 333

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 79
$wsShow.Range("F13").Value = 62

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 79
$wsAll.Range("F7").Value = 675
$wsAll.Range("F8").Value = 934
$wsAll.Range("F9").Value = 697
$wsAll.Range("F18").Value = 371
$wsAll.Range("F22").Value = 332
$wsAll.Range("F24").Value = 333
$wsAll.Range("F32").Value = 62
